$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The paragraph that used to read "Billy has some" becomes the
#    lead-in sentence for a Billy Wilder quote taken from Cameron
#    Crowe's "Conversations With Billy Wilder".
# ------------------------------------------------------------------
$lead = $d.Content
$lead.Find.Execute("Billy has some")
$lead.Text = "From Cameron Crowe's"

$space = $d.Range($lead.End, $lead.End)
$space.InsertAfter(" ")

$linkText = $d.Range($space.End, $space.End)
$linkText.InsertAfter("Conversations With Billy Wilder Advice For Screenwriters")

# Turn the text we just inserted into a real hyperlink, matching the
# existing "Conversations With Billy Wilder" link elsewhere in the
# document (same target URL).
$toLink = $d.Content
$toLink.Find.Execute("Conversations With Billy Wilder Advice For Screenwriters")
$d.Hyperlinks.Add($toLink, "http://www.listsofnote.com/2012/03/advice-from-billy-wilder.html") | Out-Null

# ------------------------------------------------------------------
# 2. Add the quote itself plus its attribution as two new BlockText
#    paragraphs right after the lead-in paragraph.
# ------------------------------------------------------------------
$afterLink = $d.Content
$afterLink.Find.Execute("Conversations With Billy Wilder Advice For Screenwriters")
$afterLink.Collapse(0)
$afterLink.InsertAfter("`rAn actor entering through the door, you've got nothing. But if he enters through the window, you've got a situation.`r--Billy Wilder")

$quotePara = $d.Content
$quotePara.Find.Execute("An actor entering through the door")
$quotePara.Paragraphs(1).Style = "Block Text"

$attribPara = $d.Content
$attribPara.Find.Execute("--Billy Wilder")
$attribPara.Paragraphs(1).Style = "Block Text"
